# Fixed diagonal movement logic
# Move the corresponding ScrumBoard cards from "In Progress"/"To verify"
# into "Done" now that the diagonal-movement work has actually landed.

$wb = $excel.ActiveWorkbook
$wsUserStories = $wb.Worksheets.Item("UserStories")
$wsScrumBoard  = $wb.Worksheets.Item("ScrumBoard")
$wsCharts      = $wb.Worksheets.Item("Charts")

# --- ScrumBoard: move cards from column C (In Progress) / D (To verify)
#     into column E (Done) ---

# Row 4: W1 (owner Lesli) In Progress -> Done
$wsScrumBoard.Cells.Item(4, 3).Clear()
$wsScrumBoard.Cells.Item(4, 5).Value2 = "W1"

# Row 5: W2 (owner Hannah) In Progress -> Done
$wsScrumBoard.Cells.Item(5, 3).Clear()
$wsScrumBoard.Cells.Item(5, 5).Value2 = "W2"

# Row 6: the card tracked here is now W2 (owner Hannah) and lands in Done
$wsScrumBoard.Cells.Item(6, 3).Clear()
$wsScrumBoard.Cells.Item(6, 5).Value2 = "W2"

# Row 17: C7 (owner Brandon) To verify -> Done
$wsScrumBoard.Cells.Item(17, 4).Clear()
$wsScrumBoard.Cells.Item(17, 5).Value2 = "C7"

# Row 20: B1 (owner Brandon) To verify -> Done
$wsScrumBoard.Cells.Item(20, 4).Clear()
$wsScrumBoard.Cells.Item(20, 5).Value2 = "B1"

# --- Selections / active sheet bookkeeping ---

$wsUserStories.Activate() | Out-Null
$wsUserStories.Range("E10").Select() | Out-Null

$wsScrumBoard.Activate() | Out-Null
$wsScrumBoard.Range("D20").Select() | Out-Null

$wsCharts.Activate() | Out-Null
$wsCharts.Range("D11").Select() | Out-Null
